$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row for "H 72" (row 2) was removed from the data set, causing every
# subsequent row to shift up by one. Deleting the entire row 2 reproduces
# this: Excel moves rows 3:63 up to become rows 2:62 and updates the
# worksheet's used range/dimension accordingly.
$ws.Rows.Item(2).Delete()
